$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.740.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.35%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.700.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.35%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.46%  '

$ws.Range("E6").Value = '  +0.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3944'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.52%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4043'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.521'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.49'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08875'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.61%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.476'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.129'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001325'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.705.56'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '99.75'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07060'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.077'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.003'
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.77%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.727.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.210'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.372'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.67%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.823'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +18.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.09%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.177'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.785'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.81%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08945'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.078'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.70%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.991'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.88%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2757'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '14.56'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02789'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09177'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.463'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7726'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7211'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.75%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.571'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.220'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.31%  '

$ws.Range("B47").Value = 'Flow'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.363'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.43%  '

$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.30%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.62%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '91.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07996'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.82%  '
